$d = $word.ActiveDocument
$deg = [char]0x00B0   # "°"

# ------------------------------------------------------------------------
# The "VISTO il Regolamento di Organizzazione ... DPCNR n. 119 ... ;"
# paragraph needs its decree reference updated to the newer one. The run
# ends up split into three pieces:
#   1) "...DPCNR n. "                                  (untouched prefix)
#   2) "144 Prot. n. 521963 del 19 dicembre 2025, in
#       vigore dal 1° gennaio 2026"                     (new decree info)
#   3) ";"                                               (kept, re-typed)
# ------------------------------------------------------------------------

# 1) Locate the unchanged prefix; its end is where the replacement begins.
$prefixRange = $d.Content
$prefixRange.Find.Execute(
    "il Regolamento di Organizzazione e Funzionamento del Consiglio Nazionale delle Ricerche - DPCNR n. "
) | Out-Null
$splitPoint = $prefixRange.End

# 2) Locate the old decree number/protocol/date text (stop right before the
#    trailing semicolon so the semicolon can be handled separately).
$oldDetails = $d.Content
$oldDetails.Start = $splitPoint
$oldDetails.Find.Execute(
    "119 prot. n. 241776 del 10 luglio 2024, entrato in vigore dal 1" + $deg + " agosto 2024"
) | Out-Null

# Briefly toggling a character attribute around the text assignment makes
# the engine materialize the new text in its own run instead of merging it
# back into the original one - the same run-splitting behaviour Word shows
# when freshly typed text replaces a selection.
$newDetails = "144 Prot. n. 521963 del 19 dicembre 2025, in vigore dal 1" + $deg + " gennaio 2026"
$oldDetails.Font.Bold = 1
$oldDetails.Text = $newDetails
$detailsEnd = $oldDetails.Start + $newDetails.Length
$newDetailsRange = $d.Range($oldDetails.Start, $detailsEnd)
$newDetailsRange.Font.Bold = 0

# 3) Re-type the trailing semicolon too, so it becomes its own run rather
#    than remaining part of the (now differently-worded) original run.
$semicolon = $d.Range($detailsEnd, $detailsEnd + 1)
$semicolon.Font.Bold = 1
$semicolon.Text = "X"
$d.Range($detailsEnd, $detailsEnd + 1).Text = ";"
$d.Range($detailsEnd, $detailsEnd + 1).Font.Bold = 0
